$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $origStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = $origStyle
}

Set-TextValue $ws.Range('D2') '29.782.79'
Set-TextValue $ws.Range('E2') '  -0.39%  '
Set-TextValue $ws.Range('D3') '1.888.28'
Set-TextValue $ws.Range('E3') '  -0.78%  '
Set-TextValue $ws.Range('E4') '  +0.09%  '
Set-TextValue $ws.Range('D5') '0.7931'
Set-TextValue $ws.Range('E5') '  -1.26%  '
Set-TextValue $ws.Range('D6') '241.67'
Set-TextValue $ws.Range('E6') '  +0.41%  '
Set-TextValue $ws.Range('D7') '1.001'
Set-TextValue $ws.Range('E7') '  +0.09%  '
Set-TextValue $ws.Range('D8') '0.3170'
Set-TextValue $ws.Range('E8') '  +2.03%  '
Set-TextValue $ws.Range('D9') '25.45'
Set-TextValue $ws.Range('E9') '  -3.27%  '
Set-TextValue $ws.Range('E10') '  +0.66%  '
Set-TextValue $ws.Range('D11') '0.08052'
Set-TextValue $ws.Range('E11') '  +0.79%  '
Set-TextValue $ws.Range('D12') '0.7661'
Set-TextValue $ws.Range('E12') '  +3.90%  '
Set-TextValue $ws.Range('D13') '1.920.42'
Set-TextValue $ws.Range('E13') '  +1.42%  '
Set-TextValue $ws.Range('D14') '5.294'
Set-TextValue $ws.Range('E14') '  +2.73%  '
Set-TextValue $ws.Range('D15') '92.09'
Set-TextValue $ws.Range('E15') '  -0.07%  '
Set-TextValue $ws.Range('D16') '29.799.39'
Set-TextValue $ws.Range('E16') '  -0.38%  '
Set-TextValue $ws.Range('D17') '13.79'
Set-TextValue $ws.Range('E17') '  -0.87%  '
Set-TextValue $ws.Range('D18') '5.929'
Set-TextValue $ws.Range('E18') '  +1.56%  '
Set-TextValue $ws.Range('D19') '242.79'
Set-TextValue $ws.Range('E19') '  -0.43%  '
Set-TextValue $ws.Range('E20') '  -0.90%  '
Set-TextValue $ws.Range('D22') '2.154.31'
Set-TextValue $ws.Range('E22') '  -0.71%  '
Set-TextValue $ws.Range('D23') '8.089'
Set-TextValue $ws.Range('E23') '  +17.56%  '
Set-TextValue $ws.Range('D24') '1.001'
Set-TextValue $ws.Range('E24') '  +0.10%  '
Set-TextValue $ws.Range('D25') '0.1621'
Set-TextValue $ws.Range('E25') '  +11.70%  '
Set-TextValue $ws.Range('D26') '9.285'
Set-TextValue $ws.Range('E26') '  +1.31%  '
Set-TextValue $ws.Range('D27') '163.80'
Set-TextValue $ws.Range('E27') '  -2.26%  '
Set-TextValue $ws.Range('E28') '  -0.92%  '
Set-TextValue $ws.Range('D29') '2.054'
Set-TextValue $ws.Range('E29') '  -0.21%  '
Set-TextValue $ws.Range('D30') '1.371'
Set-TextValue $ws.Range('E30') '  +1.23%  '
Set-TextValue $ws.Range('D31') '1.537'
Set-TextValue $ws.Range('E31') '  +1.73%  '
Set-TextValue $ws.Range('D32') '4.440'
Set-TextValue $ws.Range('E32') '  +3.98%  '
Set-TextValue $ws.Range('D33') '0.05639'
Set-TextValue $ws.Range('E33') '  +2.39%  '
Set-TextValue $ws.Range('D34') '4.086'
Set-TextValue $ws.Range('E34') '  +0.98%  '
Set-TextValue $ws.Range('D35') '1.262'
Set-TextValue $ws.Range('E35') '  +0.25%  '
Set-TextValue $ws.Range('D36') '0.7356'
Set-TextValue $ws.Range('E36') '  +0.81%  '
Set-TextValue $ws.Range('D37') '1.001'
Set-TextValue $ws.Range('E37') '  +0.22%  '
Set-TextValue $ws.Range('E38') '  -0.15%  '
Set-TextValue $ws.Range('D39') '0.01921'
Set-TextValue $ws.Range('E39') '  +0.23%  '
Set-TextValue $ws.Range('D40') '2.769'
Set-TextValue $ws.Range('E40') '  -0.51%  '
Set-TextValue $ws.Range('D41') '0.4414'
Set-TextValue $ws.Range('D42') '72.10'
Set-TextValue $ws.Range('E42') '  +0.07%  '
Set-TextValue $ws.Range('D43') '5.823'
Set-TextValue $ws.Range('E43') '  -2.25%  '
Set-TextValue $ws.Range('E44') '  +0.09%  '
Set-TextValue $ws.Range('D45') '0.8408'
Set-TextValue $ws.Range('E45') '  +0.64%  '
Set-TextValue $ws.Range('D46') '1.028.55'
Set-TextValue $ws.Range('E46') '  +5.05%  '
Set-TextValue $ws.Range('B47') 'Quant'
Set-TextValue $ws.Range('C47') 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
Set-TextValue $ws.Range('D47') '101.95'
Set-TextValue $ws.Range('E47') '  +1.40%  '
Set-TextValue $ws.Range('B48') 'EnergySwap'
Set-TextValue $ws.Range('C48') 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
Set-TextValue $ws.Range('D48') '9.963'
Set-TextValue $ws.Range('E48') '  +3.04%  '
Set-TextValue $ws.Range('B49') 'RenderToken'
Set-TextValue $ws.Range('C49') 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue $ws.Range('D49') '1.862'
Set-TextValue $ws.Range('E49') '  -1.09%  '
Set-TextValue $ws.Range('D50') '7.431'
Set-TextValue $ws.Range('E50') '  -1.22%  '
Set-TextValue $ws.Range('D51') '2.040.85'
Set-TextValue $ws.Range('E51') '  -1.19%  '
